$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while we assign numeric-looking strings,
# so Excel does not silently convert them to Double values, then restore
# the original (default) style so the cell style/index is unaffected.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.441.71"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.105.38"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "334.21"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "0.5231"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "0.4556"
$ws.Range("E8").Value = "  +5.79%  "
$ws.Range("D9").Value = "53.31"
$ws.Range("E9").Value = "  +15.59%  "
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "24.17"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "2.104.37"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "8.076"
$ws.Range("E15").Value = "  +5.28%  "
$ws.Range("D16").Value = "96.67"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "0.06652"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "19.24"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "6.344"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "30.502.03"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").Value = "12.52"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Value = "2.349"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "2.353.85"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "22.27"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "162.81"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "133.24"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "1.214"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "1.670"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "6.401"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").Value = "3.938"
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("D36").Value = "10.42"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").Value = "5.766"
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "0.06855"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("D40").Value = "0.2300"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("D41").Value = "12.73"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "0.6896"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "14.05"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "2.318"
$ws.Range("E45").Value = "  +5.32%  "
$ws.Range("D46").Value = "0.6380"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "3.661"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "0.00000000348"
$ws.Range("E49").Value = "  +20.99%  "
$ws.Range("D50").Value = "83.49"
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "1.206"
$ws.Range("E51").Value = "  +1.70%  "

$ws.Range("D2:D51").Style = "Normal"

